$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.371.53"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "2.510.06"
$ws.Range("E3").Value = "  -4.53%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.88"
$c.ClearFormats()
$ws.Range("E5").Value = "  -2.17%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "172.62"
$c.ClearFormats()
$ws.Range("E6").Value = "  +2.40%  "
$ws.Range("E7").Value = "  +0.10%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.522"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "2.509.65"
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("E11").Value = "  -0.35%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.349"
$c.ClearFormats()
$ws.Range("E12").Value = "  -4.16%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.10"
$c.ClearFormats()
$ws.Range("E13").Value = "  -2.30%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "26.51"
$c.ClearFormats()
$ws.Range("E14").Value = "  -4.17%  "
$ws.Range("D15").Value = "2.966.32"
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("E16").Value = "  -3.93%  "
$ws.Range("D17").Value = "66.319.31"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D18").Value = "2.508.24"
$ws.Range("E18").Value = "  -4.87%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.23"
$c.ClearFormats()
$ws.Range("E19").Value = "  -6.51%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.67"
$c.ClearFormats()
$ws.Range("E20").Value = "  -4.88%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "347.54"
$c.ClearFormats()
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("E22").Value = "  -2.85%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.58"
$c.ClearFormats()
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("E25").Value = "  +0.02%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "69.48"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.45%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.84"
$c.ClearFormats()
$ws.Range("E27").Value = "  -4.41%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "2.636.52"
$ws.Range("E29").Value = "  -4.48%  "
$ws.Range("D30").Value = "0.0₃0974"
$ws.Range("E30").Value = "  -3.55%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "527.65"
$c.ClearFormats()
$ws.Range("E31").Value = "  -3.75%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.06"
$c.ClearFormats()
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -3.47%  "
$ws.Range("E36").Value = "  -0.03%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "157.51"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  +0.29%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.354"
$c.ClearFormats()
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  -0.06%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.48"
$c.ClearFormats()
$ws.Range("E45").Value = "  +1.96%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "39.56"
$c.ClearFormats()
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "147.12"
$c.ClearFormats()
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.556"
$c.ClearFormats()
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "3.66"
$c.ClearFormats()
$ws.Range("E50").Value = "  -3.56%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0268"
$ws.Range("E51").Value = "  -9.96%  "
